$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.270.13'
$ws.Range('E2').Value = '  -0.92%  '
$ws.Range('D3').Value = '1.866.50'
$ws.Range('E3').Value = '  -1.86%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.11'
$ws.Range('E5').Value = '  -1.61%  '
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4674'
$ws.Range('E7').Value = '  -1.25%  '
$ws.Range('E8').Value = '  -0.72%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06549'
$ws.Range('E9').Value = '  -1.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.77'
$ws.Range('E10').Value = '  +5.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07873'
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '97.03'
$ws.Range('E12').Value = '  -3.39%  '
$ws.Range('D13').Value = '1.867.00'
$ws.Range('E13').Value = '  -2.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.142'
$ws.Range('E14').Value = '  -0.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6753'
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '280.95'
$ws.Range('E16').Value = '  -0.88%  '
$ws.Range('D17').Value = '30.271.59'
$ws.Range('E17').Value = '  -1.01%  '
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.497'
$ws.Range('E19').Value = '  +1.65%  '
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('D21').Value = '2.116.10'
$ws.Range('E21').Value = '  -1.85%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.000007260'
$ws.Range('E22').Value = '  -2.86%  '
$ws.Range('E23').Value = '  -0.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.181'
$ws.Range('E24').Value = '  -1.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.302'
$ws.Range('E25').Value = '  -0.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.12'
$ws.Range('E26').Value = '  -1.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.11'
$ws.Range('E27').Value = '  -0.86%  '
$ws.Range('E28').Value = '  -5.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.350'
$ws.Range('E29').Value = '  -2.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09663'
$ws.Range('E30').Value = '  -2.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.415'
$ws.Range('E31').Value = '  -1.93%  '
$ws.Range('E32').Value = '  -2.71%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.103'
$ws.Range('E33').Value = '  -3.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04705'
$ws.Range('E34').Value = '  -1.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7039'
$ws.Range('E35').Value = '  -2.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.103'
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01860'
$ws.Range('E38').Value = '  -2.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.371'
$ws.Range('E39').Value = '  -4.81%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.531'
$ws.Range('E40').Value = '  -1.74%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '73.42'
$ws.Range('E41').Value = '  -0.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.940'
$ws.Range('E42').Value = '  -2.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8478'
$ws.Range('E43').Value = '  -2.79%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '104.32'
$ws.Range('E44').Value = '  -0.45%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4177'
$ws.Range('E45').Value = '  -2.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9995'
$ws.Range('E46').Value = '  -0.22%  '
$ws.Range('E47').Value = '  -2.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.244'
$ws.Range('E48').Value = '  +0.42%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '936.80'
$ws.Range('E49').Value = '  -4.98%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.12'
$ws.Range('E50').Value = '  -1.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1133'
$ws.Range('E51').Value = '  -4.27%  '
